# Update Costa Rica Primera Division odds sheet
# - swap the two teams "Puntarenas" / "AD Guanacasteca" (shared string re-labelling)
# - re-order several match rows that share an identical kickoff timestamp
# - append the newest scraped match (row 323)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-order rows that were re-fetched in a different sequence for matches
#    that share the exact same Date/time value. Column A (running id),
#    C (Div), D (Div Original Name) and E (Date) stay put - only the
#    match-specific data in B and F:AC move between the two rows.
# ---------------------------------------------------------------------------
$swapPairs = @(
    @(8,9),
    @(81,82),
    @(104,105),
    @(124,125),
    @(178,179),
    @(231,232),
    @(250,251),
    @(269,271)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value = $v2
    $range2.Value = $v1
}

# Rows 130-133 were reshuffled as a 4-way rotation instead of a simple swap
$v130 = $ws.Range("B130:AC130").Value2
$v131 = $ws.Range("B131:AC131").Value2
$v132 = $ws.Range("B132:AC132").Value2
$v133 = $ws.Range("B133:AC133").Value2
$ws.Range("B130:AC130").Value = $v131
$ws.Range("B131:AC131").Value = $v132
$ws.Range("B132:AC132").Value = $v133
$ws.Range("B133:AC133").Value = $v130

# ---------------------------------------------------------------------------
# 2) Relabel the two teams whose names were swapped in this refresh.
#    Use a temporary placeholder so the two find/replace passes cannot
#    collide with each other.
# ---------------------------------------------------------------------------
$ws.Cells.Replace("Puntarenas", "__TEAM_SWAP_TMP__", 1)
$ws.Cells.Replace("AD Guanacasteca", "Puntarenas", 1)
$ws.Cells.Replace("__TEAM_SWAP_TMP__", "AD Guanacasteca", 1)

# ---------------------------------------------------------------------------
# 3) Append the newly scraped match as row 323
# ---------------------------------------------------------------------------
$newRow = 323
$ws.Range("A$newRow").Value = 321
$ws.Range("B$newRow").Value = 7623912
$ws.Range("C$newRow").Value = "Costa Rica Primera Division"
$ws.Range("D$newRow").Value = "Costa Rica Primera Division"
$ws.Range("E$newRow").Value = 45340.875
$ws.Range("F$newRow").Value = "Herediano"
$ws.Range("G$newRow").Value = "AD Guanacasteca"
$ws.Range("H$newRow").Value = 2
$ws.Range("I$newRow").Value = 1
$ws.Range("J$newRow").Value = "H"
$ws.Range("K$newRow").Value = 1.8
$ws.Range("L$newRow").Value = 3.5
$ws.Range("M$newRow").Value = 3.75
$ws.Range("N$newRow").Value = 1.5
$ws.Range("O$newRow").Value = 3.8
$ws.Range("P$newRow").Value = 5.5
$ws.Range("Q$newRow").Value = -1
$ws.Range("R$newRow").Value = 1.85
$ws.Range("S$newRow").Value = 1.95
$ws.Range("T$newRow").Value = 2.75
$ws.Range("U$newRow").Value = 2
$ws.Range("V$newRow").Value = 1.8
$ws.Range("W$newRow").Value = 0.5
$ws.Range("X$newRow").Value = -1
$ws.Range("Y$newRow").Value = -1
$ws.Range("Z$newRow").Value = 0
$ws.Range("AA$newRow").Value = 0
$ws.Range("AB$newRow").Value = 0.5
$ws.Range("AC$newRow").Value = -0.5

$ws.Range("A$newRow").Font.Bold = $ws.Range("A322").Font.Bold
